# "Draft final interim report 2" — updates to the `discounting` sheet:
#   - AI5 / AI6 / AI12 / AI13: 175 -> 145 (Number of introductions over 10 years)
#   - R11: 2200 -> 130 (Avg total cases in the popn over 10 years ...)
#   - sheetView: topLeftCell AH1 -> AI1, selection BA4 -> BC11

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("discounting")

$ws.Range("AI5").Value = 145
$ws.Range("AI6").Value = 145
$ws.Range("R11").Value = 130
$ws.Range("AI12").Value = 145
$ws.Range("AI13").Value = 145

$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 35
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("BC11").Select()
